# Escandallo Aveiro P6.76 — add M3 nuts to the power supply (fuente) section
# and adjust the "Cable plano 20 vías" items' lengths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the two "Cable plano 20 vías" descriptions --------------
$ws.Range("A10").Value = "Cable plano 20 vías 40 cm H-H"
$ws.Range("A11").Value = "Cable plano 20 vías 90 cm H-H"

# --- 2. Append the new "FUENTE DE ALIMENTACIÓN" rows --------------------
# Row 45: Torreta L=30 H-M  x4
# Row 46: Tuerca M3         x10
$ws.Range("A44:B44").Copy()
$ws.Range("A45:B46").PasteSpecial(-4122)

$ws.Range("A45").Value = "Torreta L=30 H-M"
$ws.Range("B45").Value = 4
$ws.Range("A46").Value = "Tuerca M3 "
$ws.Range("B46").Value = 10

# --- 3. Print area now covers the extended sheet + the diagram image ----
$ws.PageSetup.PrintArea = '$A$1:$E$101'
$ws.PageSetup.Zoom = 58

# --- 4. Reposition / resize the title text box and diagram image --------
$title = $ws.Shapes.Item("Marco de texto 1")
$title.Top = 13.373622047244094
$title.Left = 234.8787401574803
$title.Width = 453.6212598425197
$title.Height = 35.36456692913386

$diagram = $ws.Shapes.Item("Imagen 5")
$diagram.Top = 601.1785826771653
$diagram.Left = 2.5000787401574804
$diagram.Width = 622.1427559055118
$diagram.Height = 699.3066141732284

# --- 5. View state: zoom + selection -------------------------------------
$excel.ActiveWindow.Zoom = 70
$ws.Range("D49").Select()
